# Hjemme passive updated meanEMG legmaxROM
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 1 header values updated
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 values updated; B2 cleared entirely
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = 0.71303677533663179
$ws.Range("D2").Value = 0.425890916456372
$ws.Range("E2").Value = 1.0810840941366469

# Row 3 values updated
$ws.Range("B3").Value = 1.5369245927528954
$ws.Range("C3").Value = 1.9833558962570397
$ws.Range("D3").Value = 6.9989719429577564
$ws.Range("E3").Value = 5.0498057792531421

# Update selection to match new active range
$ws.Range("B1:E3").Select()
